$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad date) from row 2 through row 56: 45212 -> 45221
for ($r = 2; $r -le 56; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45212) {
        $cell.Value2 = 45221
    }
}
